$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","H","I","L","N")

$data = @(
    @(1.699955922942479,0.2285688705859457,0.0781804742505301,0.05614855279578546,1.29221744435344,1.164334417649656,1.024765759105783,0.2596427119920719,1.348250203937496),
    @(1.572168292037361,0.1991593398948339,0.07094743015032634,0.05646580474595986,1.265474754307547,1.158506721433213,1.0222892936915,0.2492817122659261,1.367548845873078),
    @(1.494441646628331,0.1810649040776866,0.06654695993825044,0.05667718952738898,1.250061537481571,1.155649144840595,1.021452894491482,0.2430745017239957,1.379976849622762),
    @(1.46295195584122,0.1736816593498247,0.06476381703531331,0.05676750780481399,1.244032242158639,1.154665336481571,1.021283465296307,0.2405837446554102,1.385186595943261),
    @(1.45773425212235,0.1724550864731214,0.06446833261733786,0.05678275754138618,1.243046242257165,1.154512872749052,1.021265666213587,0.2401724923401929,1.386060434907026),
    @(1.494016218924287,0.1809653704376331,0.06652287119564448,0.05667839067031544,1.24997920706727,1.155635146021751,1.02144991631441,0.2430407538285664,1.38004652237737),
    @(1.6557418135136,0.2184358985144286,0.07567799923475604,0.05625450105786367,1.282786677944898,1.162175126916139,1.023769493600263,0.2560381022171896,1.354784146659205),
    @(1.978757495564707,0.2916421217224183,0.09396084055697429,0.05555466253301411,1.35518338766434,1.180744741693871,1.033776831652752,0.2827585433657163,1.309845746343682),
    @(2.219734940859951,0.3452917259901938,0.1076054907163382,0.05512029256560869,1.413395134398428,1.197930495944689,1.044501369716642,0.303153585871641,1.279647673503747),
    @(2.330174917434476,0.3696761919278515,0.1138614121524029,0.05493995108754746,1.440991441880499,1.20652707251233,1.050122310712382,0.3126004037800669,1.266524643963094),
    @(2.372114303767319,0.3789074104468568,0.1162375757570544,0.05487413691990106,1.451603655706833,1.209895034249712,1.05235832147433,0.3162021508728259,1.261643913523844),
    @(2.363076642580438,0.3769194142326455,0.1157255045867487,0.05488820106723136,1.449310889529158,1.209164663961644,1.05187196434737,0.3154253608807949,1.262691116585962),
    @(2.333622924731685,0.3704357007831049,0.1140567557609842,0.05493448689212688,1.441861256842088,1.206801895828988,1.050304110281743,0.3128962310182857,1.266121324095148),
    @(2.315597076424069,0.3664639051600034,0.1130355386980284,0.05496316078536534,1.437319301050138,1.205369318635604,1.049357773914366,0.3113502532106907,1.268233984025983),
    @(2.212533935118586,0.3436977526739042,0.1071976483432309,0.05513242518000627,1.411614240987689,1.197384414053118,1.044149028057632,0.3025396268133136,1.28051768398047),
    @(2.149517937839335,0.329726337266635,0.103628922752705,0.05524068003578098,1.396131916909383,1.192685834712819,1.041144275595968,0.297177978273524,1.288210901369451),
    @(2.113349776200323,0.3216883486981601,0.1015808799736249,0.05530456978467058,1.387331787257665,1.190056590818102,1.03948584496468,0.2941099994625205,1.292693628880924),
    @(2.101117063047013,0.3189664691678331,0.100888232325957,0.05532648087744096,1.384370179081458,1.189178935902987,1.038936298638419,0.2930739611533113,1.294221318693417),
    @(2.15621812822684,0.3312138225898025,0.1040083431597481,0.05522898803012843,1.397769168511701,1.193178420275359,1.041456904836878,0.2977470883669753,1.287385961465529),
    @(2.342270984417553,0.3723401935083643,0.1145467116224808,0.05492082442420188,1.444044981927703,1.207492836274525,1.05076170403602,0.3136384332921267,1.26511137952137),
    @(2.464556293505325,0.399203387528928,0.1214760771485288,0.05473385907525863,1.475234485994577,1.217504906066125,1.057469767709435,0.3241668801018278,1.251070692661607),
    @(2.399227099043344,0.3848672679652623,0.1177738575609624,0.05483232631264467,1.458500978553133,1.21210096421612,1.053831944959867,0.3185345613842543,1.258517044300673),
    @(2.153188783512348,0.3305413480636901,0.1038367957913096,0.05523426883828009,1.397028653008419,1.192955498012964,1.041315350151827,0.2974897485271697,1.287758730561808),
    @(1.890736575784786,0.2718639136744798,0.08897839490062154,0.05572995144762682,1.33472408446579,1.175101999150172,1.03048063298133,0.2753967533701456,1.321508858949683)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowNum = $r + 2
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $data[$r][$c]
    }
}
